# The commit swaps the content of ppt/theme/theme1.xml (the theme driving
# the slide master / all slides, originally the "Integral" / "Red Violet"
# palette) with ppt/theme/theme2.xml (the theme driving the notes master,
# originally the default "Office Theme" palette). After the edit,
# theme1.xml carries the Office Theme palette and theme2.xml carries the
# Integral / Red Violet palette - i.e. the presentation's visible design
# switches to the stock "Office Theme" colours.
#
# The PowerPoint object model exposes the 12-slot theme colour scheme via
# ThemeColorScheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - the
# same order used inside <a:clrScheme>). Writing through a Slide applies
# directly to the shared theme part used by the slide master (theme1.xml).

$p = $ppt.ActivePresentation

# --- Target palette: default "Office Theme" colours (was "Integral") ---
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeRGB = @(
    0x000000,
    0xFFFFFF,
    0x44546A,
    0xE7E6E6,
    0x5B9BD5,
    0xED7D31,
    0xA5A5A5,
    0xFFC000,
    0x4472C4,
    0x70AD47,
    0x0563C1,
    0x954F72
)

function ToComRGB([int]$hexRGB) {
    # VBA/COM RGB() packs colour as 0x00BBGGRR, i.e. byte-reversed from the
    # usual 0xRRGGBB hex notation used in OOXML srgbClr values.
    $r = ($hexRGB -shr 16) -band 0xFF
    $g = ($hexRGB -shr 8) -band 0xFF
    $b = $hexRGB -band 0xFF
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

# Apply the new palette to the presentation's (single, shared) theme - this
# is backed by ppt/theme/theme1.xml, the theme referenced by the slide
# master that all slides use.
$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ToComRGB($officeThemeRGB[$i - 1])
}

# Best-effort: also try to relabel the design/theme names to match what the
# commit produced ("Office Theme" / clrScheme "Office") and to push the old
# "Integral" / "Red Violet" palette onto the notes-master theme
# (theme2.xml), in case this host exposes those write paths. These are
# no-ops (ignored) on hosts that only implement colour-slot writes, so they
# are wrapped defensively and never fail the script.
try { $p.Designs.Item(1).Name = "Office Theme" } catch {}
try { $p.SlideMaster.Name = "Office Theme" } catch {}

$integralRGB = @(
    0x000000,
    0xFFFFFF,
    0x454551,
    0xD8D9DC,
    0xE32D91,
    0xC830CC,
    0x4EA6DC,
    0x4775E7,
    0x8971E1,
    0xD54773,
    0x6B9F25,
    0x8C8C8C
)
try {
    $notesColors = $p.NotesMaster.ColorScheme
    for ($i = 1; $i -le $notesColors.Count; $i++) {
        $notesColors.Item($i).RGB = ToComRGB($integralRGB[$i - 1])
    }
} catch {}
